$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was row 5)
$ws.Range("D2").Value = 44330
$ws.Range("K2").Value = 'Mankaki'
$ws.Range("L2").Value = 'Primera'
$ws.Range("M2").Value = 55
$ws.Range("N2").Value = 18000
$ws.Range("O2").Value = 18000
$ws.Range("P2").Value = 18000
$ws.Range("Q2").Value = '$/bandeja 15 kilos empedrada'
$ws.Range("R2").Value = 'Región de O''Higgins'
$ws.Range("S2").Value = 1200
$ws.Range("T2").Value = 15

# Row 3 (was row 50)
$ws.Range("D3").Value = 44315
$ws.Range("K3").Value = 'Fuyu'
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 85
$ws.Range("N3").Value = 16000
$ws.Range("O3").Value = 17000
$ws.Range("P3").Value = 16471
$ws.Range("Q3").Value = '$/bandeja 15 kilos granel'
$ws.Range("R3").Value = 'Región de O''Higgins'
$ws.Range("S3").Value = 1098
$ws.Range("T3").Value = 15

# Row 4 (was row 3)
$ws.Range("D4").Value = 44698
$ws.Range("K4").Value = 'Mankaki'
$ws.Range("L4").Value = 'Especial'
$ws.Range("M4").Value = 30
$ws.Range("N4").Value = 18000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 18000
$ws.Range("Q4").Value = '$/bandeja 15 kilos granel'
$ws.Range("R4").Value = 'Región de O''Higgins'
$ws.Range("S4").Value = 1200
$ws.Range("T4").Value = 15

# Row 5 (was row 11)
$ws.Range("D5").Value = 44729
$ws.Range("K5").Value = 'Mankaki'
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 30
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 20000
$ws.Range("Q5").Value = '$/bandeja 15 kilos granel'
$ws.Range("R5").Value = 'Región de O''Higgins'
$ws.Range("S5").Value = 1333
$ws.Range("T5").Value = 15

# Row 6 (was row 45)
$ws.Range("D6").Value = 44349
$ws.Range("K6").Value = 'Fuyu'
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 14000
$ws.Range("O6").Value = 14000
$ws.Range("P6").Value = 14000
$ws.Range("Q6").Value = '$/bandeja 15 kilos granel'
$ws.Range("R6").Value = 'Región de O''Higgins'
$ws.Range("S6").Value = 933
$ws.Range("T6").Value = 15

# Row 7 (was row 46)
$ws.Range("D7").Value = 44349
$ws.Range("K7").Value = 'Mankaki'
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 80
$ws.Range("N7").Value = 14000
$ws.Range("O7").Value = 14000
$ws.Range("P7").Value = 14000
$ws.Range("Q7").Value = '$/bandeja 15 kilos granel'
$ws.Range("R7").Value = 'Región de O''Higgins'
$ws.Range("S7").Value = 933
$ws.Range("T7").Value = 15

# Row 8 (was row 9)
$ws.Range("D8").Value = 44316
$ws.Range("K8").Value = 'Fuyu'
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 110
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 16000
$ws.Range("P8").Value = 15409
$ws.Range("Q8").Value = '$/bandeja 15 kilos granel'
$ws.Range("R8").Value = 'Provincia de Limarí'
$ws.Range("S8").Value = 1027
$ws.Range("T8").Value = 15

# Row 9 (was row 20)
$ws.Range("D9").Value = 44305
$ws.Range("K9").Value = 'Fuyu'
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 80
$ws.Range("N9").Value = 17000
$ws.Range("O9").Value = 17000
$ws.Range("P9").Value = 17000
$ws.Range("Q9").Value = '$/bandeja 15 kilos granel'
$ws.Range("R9").Value = 'Provincia de Limarí'
$ws.Range("S9").Value = 1133
$ws.Range("T9").Value = 15

# Row 10 (was row 21)
$ws.Range("D10").Value = 44305
$ws.Range("K10").Value = 'Fuyu'
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 120
$ws.Range("N10").Value = 10000
$ws.Range("O10").Value = 11000
$ws.Range("P10").Value = 10583
$ws.Range("Q10").Value = '$/caja 10 kilos granel'
$ws.Range("R10").Value = 'Provincia de Limarí'
$ws.Range("S10").Value = 10583
$ws.Range("T10").Value = 1

# Row 11 (was row 22)
$ws.Range("D11").Value = 44305
$ws.Range("K11").Value = 'Mankaki'
$ws.Range("L11").Value = 'Primera'
$ws.Range("M11").Value = 80
$ws.Range("N11").Value = 10000
$ws.Range("O11").Value = 11000
$ws.Range("P11").Value = 10500
$ws.Range("Q11").Value = '$/caja 10 kilos granel'
$ws.Range("R11").Value = 'Provincia de Limarí'
$ws.Range("S11").Value = 10500
$ws.Range("T11").Value = 1

# Row 14 (was row 47)
$ws.Range("D14").Value = 44322
$ws.Range("K14").Value = 'Fuyu'
$ws.Range("L14").Value = 'Especial'
$ws.Range("M14").Value = 80
$ws.Range("N14").Value = 23000
$ws.Range("O14").Value = 23000
$ws.Range("P14").Value = 23000
$ws.Range("Q14").Value = '$/bandeja 15 kilos granel'
$ws.Range("R14").Value = 'Provincia de Limarí'
$ws.Range("S14").Value = 1533
$ws.Range("T14").Value = 15

# Row 15 (was row 16)
$ws.Range("D15").Value = 44680
$ws.Range("K15").Value = 'Fuyu'
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 35
$ws.Range("N15").Value = 15000
$ws.Range("O15").Value = 15000
$ws.Range("P15").Value = 15000
$ws.Range("Q15").Value = '$/bandeja 15 kilos'
$ws.Range("R15").Value = 'Región de O''Higgins'
$ws.Range("S15").Value = 1000
$ws.Range("T15").Value = 15

# Row 16 (was row 10)
$ws.Range("D16").Value = 44309
$ws.Range("K16").Value = 'Fuyu'
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 70
$ws.Range("N16").Value = 16000
$ws.Range("O16").Value = 17000
$ws.Range("P16").Value = 16429
$ws.Range("Q16").Value = '$/bandeja 15 kilos granel'
$ws.Range("R16").Value = 'Provincia de Limarí'
$ws.Range("S16").Value = 1095
$ws.Range("T16").Value = 15

# Row 17 (was row 33)
$ws.Range("D17").Value = 44677
$ws.Range("K17").Value = 'Fuyu'
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 55
$ws.Range("N17").Value = 15000
$ws.Range("O17").Value = 15000
$ws.Range("P17").Value = 15000
$ws.Range("Q17").Value = '$/bandeja 15 kilos'
$ws.Range("R17").Value = 'Región de O''Higgins'
$ws.Range("S17").Value = 1000
$ws.Range("T17").Value = 15

# Row 18 (was row 48)
$ws.Range("D18").Value = 44719
$ws.Range("K18").Value = 'Fuyu'
$ws.Range("L18").Value = 'Primera'
$ws.Range("M18").Value = 35
$ws.Range("N18").Value = 20000
$ws.Range("O18").Value = 20000
$ws.Range("P18").Value = 20000
$ws.Range("Q18").Value = '$/bandeja 15 kilos granel'
$ws.Range("R18").Value = 'Región de O''Higgins'
$ws.Range("S18").Value = 1333
$ws.Range("T18").Value = 15

# Row 19 (was row 49)
$ws.Range("D19").Value = 44719
$ws.Range("K19").Value = 'Mankaki'
$ws.Range("L19").Value = 'Primera'
$ws.Range("M19").Value = 35
$ws.Range("N19").Value = 20000
$ws.Range("O19").Value = 20000
$ws.Range("P19").Value = 20000
$ws.Range("Q19").Value = '$/bandeja 15 kilos granel'
$ws.Range("R19").Value = 'Región de O''Higgins'
$ws.Range("S19").Value = 1333
$ws.Range("T19").Value = 15

# Row 20 (was row 38)
$ws.Range("D20").Value = 44294
$ws.Range("K20").Value = 'Fuyu'
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 30
$ws.Range("N20").Value = 20000
$ws.Range("O20").Value = 20000
$ws.Range("P20").Value = 20000
$ws.Range("Q20").Value = '$/bandeja 15 kilos granel'
$ws.Range("R20").Value = 'Región de O''Higgins'
$ws.Range("S20").Value = 1333
$ws.Range("T20").Value = 15

# Row 21 (was row 39)
$ws.Range("D21").Value = 44294
$ws.Range("K21").Value = 'Mankaki'
$ws.Range("L21").Value = 'Primera'
$ws.Range("M21").Value = 20
$ws.Range("N21").Value = 22000
$ws.Range("O21").Value = 22000
$ws.Range("P21").Value = 22000
$ws.Range("Q21").Value = '$/bandeja 15 kilos granel'
$ws.Range("R21").Value = 'Región de O''Higgins'
$ws.Range("S21").Value = 1467
$ws.Range("T21").Value = 15

# Row 22 (was row 2)
$ws.Range("D22").Value = 44722
$ws.Range("K22").Value = 'Mankaki'
$ws.Range("L22").Value = 'Primera'
$ws.Range("M22").Value = 35
$ws.Range("N22").Value = 20000
$ws.Range("O22").Value = 20000
$ws.Range("P22").Value = 20000
$ws.Range("Q22").Value = '$/bandeja 15 kilos granel'
$ws.Range("R22").Value = 'Región de O''Higgins'
$ws.Range("S22").Value = 1333
$ws.Range("T22").Value = 15

# Row 23 (was row 32)
$ws.Range("D23").Value = 44706
$ws.Range("K23").Value = 'Mankaki'
$ws.Range("L23").Value = 'Primera'
$ws.Range("M23").Value = 45
$ws.Range("N23").Value = 18000
$ws.Range("O23").Value = 18000
$ws.Range("P23").Value = 18000
$ws.Range("Q23").Value = '$/caja 18 kilos granel'
$ws.Range("R23").Value = 'Región de O''Higgins'
$ws.Range("S23").Value = 1000
$ws.Range("T23").Value = 18

# Row 24 (was row 25)
$ws.Range("D24").Value = 44308
$ws.Range("K24").Value = 'Fuyu'
$ws.Range("L24").Value = 'Especial'
$ws.Range("M24").Value = 30
$ws.Range("N24").Value = 17000
$ws.Range("O24").Value = 17000
$ws.Range("P24").Value = 17000
$ws.Range("Q24").Value = '$/caja 15 kilos empedrada'
$ws.Range("R24").Value = 'Provincia de Limarí'
$ws.Range("S24").Value = 1133
$ws.Range("T24").Value = 15

# Row 25 (was row 14)
$ws.Range("D25").Value = 44350
$ws.Range("K25").Value = 'Fuyu'
$ws.Range("L25").Value = 'Primera'
$ws.Range("M25").Value = 90
$ws.Range("N25").Value = 15000
$ws.Range("O25").Value = 15000
$ws.Range("P25").Value = 15000
$ws.Range("Q25").Value = '$/bandeja 15 kilos granel'
$ws.Range("R25").Value = 'Región de O''Higgins'
$ws.Range("S25").Value = 1000
$ws.Range("T25").Value = 15

# Row 26 (was row 15)
$ws.Range("D26").Value = 44350
$ws.Range("K26").Value = 'Mankaki'
$ws.Range("L26").Value = 'Primera'
$ws.Range("M26").Value = 50
$ws.Range("N26").Value = 15000
$ws.Range("O26").Value = 15000
$ws.Range("P26").Value = 15000
$ws.Range("Q26").Value = '$/bandeja 15 kilos granel'
$ws.Range("R26").Value = 'Región de O''Higgins'
$ws.Range("S26").Value = 1000
$ws.Range("T26").Value = 15

# Row 27 (was row 24)
$ws.Range("D27").Value = 44307
$ws.Range("K27").Value = 'Fuyu'
$ws.Range("L27").Value = 'Primera'
$ws.Range("M27").Value = 50
$ws.Range("N27").Value = 17000
$ws.Range("O27").Value = 17000
$ws.Range("P27").Value = 17000
$ws.Range("Q27").Value = '$/bandeja 15 kilos granel'
$ws.Range("R27").Value = 'Provincia de Limarí'
$ws.Range("S27").Value = 1133
$ws.Range("T27").Value = 15

# Row 28 (was row 43)
$ws.Range("D28").Value = 44312
$ws.Range("K28").Value = 'Fuyu'
$ws.Range("L28").Value = 'Primera'
$ws.Range("M28").Value = 135
$ws.Range("N28").Value = 15000
$ws.Range("O28").Value = 16000
$ws.Range("P28").Value = 15481
$ws.Range("Q28").Value = '$/bandeja 15 kilos granel'
$ws.Range("R28").Value = 'Provincia de Limarí'
$ws.Range("S28").Value = 1032
$ws.Range("T28").Value = 15

# Row 29 (was row 23)
$ws.Range("D29").Value = 44306
$ws.Range("K29").Value = 'Hachiya'
$ws.Range("L29").Value = 'Especial'
$ws.Range("M29").Value = 30
$ws.Range("N29").Value = 20000
$ws.Range("O29").Value = 20000
$ws.Range("P29").Value = 20000
$ws.Range("Q29").Value = '$/caja 15 kilos empedrada'
$ws.Range("R29").Value = 'Región de O''Higgins'
$ws.Range("S29").Value = 1333
$ws.Range("T29").Value = 15

# Row 30 (was row 44)
$ws.Range("D30").Value = 44334
$ws.Range("K30").Value = 'Mankaki'
$ws.Range("L30").Value = 'Primera'
$ws.Range("M30").Value = 40
$ws.Range("N30").Value = 15000
$ws.Range("O30").Value = 15000
$ws.Range("P30").Value = 15000
$ws.Range("Q30").Value = '$/bandeja 15 kilos empedrada'
$ws.Range("R30").Value = 'Provincia de Limarí'
$ws.Range("S30").Value = 1000
$ws.Range("T30").Value = 15

# Row 31 (was row 52)
$ws.Range("D31").Value = 44343
$ws.Range("K31").Value = 'Mankaki'
$ws.Range("L31").Value = 'Primera'
$ws.Range("M31").Value = 65
$ws.Range("N31").Value = 16000
$ws.Range("O31").Value = 16000
$ws.Range("P31").Value = 16000
$ws.Range("Q31").Value = '$/bandeja 15 kilos granel'
$ws.Range("R31").Value = 'Región de O''Higgins'
$ws.Range("S31").Value = 1067
$ws.Range("T31").Value = 15

# Row 32 (was row 30)
$ws.Range("D32").Value = 44721
$ws.Range("K32").Value = 'Fuyu'
$ws.Range("L32").Value = 'Primera'
$ws.Range("M32").Value = 65
$ws.Range("N32").Value = 20000
$ws.Range("O32").Value = 20000
$ws.Range("P32").Value = 20000
$ws.Range("Q32").Value = '$/bandeja 15 kilos granel'
$ws.Range("R32").Value = 'Región de O''Higgins'
$ws.Range("S32").Value = 1333
$ws.Range("T32").Value = 15

# Row 33 (was row 17)
$ws.Range("D33").Value = 44351
$ws.Range("K33").Value = 'Fuyu'
$ws.Range("L33").Value = 'Primera'
$ws.Range("M33").Value = 40
$ws.Range("N33").Value = 15000
$ws.Range("O33").Value = 15000
$ws.Range("P33").Value = 15000
$ws.Range("Q33").Value = '$/bandeja 15 kilos granel'
$ws.Range("R33").Value = 'Región de O''Higgins'
$ws.Range("S33").Value = 1000
$ws.Range("T33").Value = 15

# Row 34 (was row 18)
$ws.Range("D34").Value = 44351
$ws.Range("K34").Value = 'Mankaki'
$ws.Range("L34").Value = 'Primera'
$ws.Range("M34").Value = 30
$ws.Range("N34").Value = 15000
$ws.Range("O34").Value = 15000
$ws.Range("P34").Value = 15000
$ws.Range("Q34").Value = '$/bandeja 15 kilos granel'
$ws.Range("R34").Value = 'Región de O''Higgins'
$ws.Range("S34").Value = 1000
$ws.Range("T34").Value = 15

# Row 35 (was row 6)
$ws.Range("D35").Value = 44333
$ws.Range("K35").Value = 'Fuyu'
$ws.Range("L35").Value = 'Primera'
$ws.Range("M35").Value = 100
$ws.Range("N35").Value = 12000
$ws.Range("O35").Value = 12000
$ws.Range("P35").Value = 12000
$ws.Range("Q35").Value = '$/caja 10 kilos granel'
$ws.Range("R35").Value = 'Provincia de Limarí'
$ws.Range("S35").Value = 12000
$ws.Range("T35").Value = 1

# Row 36 (was row 7)
$ws.Range("D36").Value = 44333
$ws.Range("K36").Value = 'Hachiya'
$ws.Range("L36").Value = 'Primera'
$ws.Range("M36").Value = 40
$ws.Range("N36").Value = 15000
$ws.Range("O36").Value = 15000
$ws.Range("P36").Value = 15000
$ws.Range("Q36").Value = '$/bandeja 15 kilos empedrada'
$ws.Range("R36").Value = 'Provincia de Limarí'
$ws.Range("S36").Value = 1000
$ws.Range("T36").Value = 15

# Row 37 (was row 8)
$ws.Range("D37").Value = 44333
$ws.Range("K37").Value = 'Mankaki'
$ws.Range("L37").Value = 'Primera'
$ws.Range("M37").Value = 50
$ws.Range("N37").Value = 15000
$ws.Range("O37").Value = 15000
$ws.Range("P37").Value = 15000
$ws.Range("Q37").Value = '$/bandeja 15 kilos empedrada'
$ws.Range("R37").Value = 'Provincia de Limarí'
$ws.Range("S37").Value = 1000
$ws.Range("T37").Value = 15

# Row 38 (was row 31)
$ws.Range("D38").Value = 44678
$ws.Range("K38").Value = 'Fuyu'
$ws.Range("L38").Value = 'Primera'
$ws.Range("M38").Value = 55
$ws.Range("N38").Value = 14000
$ws.Range("O38").Value = 15000
$ws.Range("P38").Value = 14636
$ws.Range("Q38").Value = '$/bandeja 15 kilos'
$ws.Range("R38").Value = 'Región de O''Higgins'
$ws.Range("S38").Value = 976
$ws.Range("T38").Value = 15

# Row 39 (was row 4)
$ws.Range("D39").Value = 44736
$ws.Range("K39").Value = 'Fuyu'
$ws.Range("L39").Value = 'Primera'
$ws.Range("M39").Value = 250
$ws.Range("N39").Value = 19000
$ws.Range("O39").Value = 20000
$ws.Range("P39").Value = 19400
$ws.Range("Q39").Value = '$/bandeja 15 kilos granel'
$ws.Range("R39").Value = 'Región de O''Higgins'
$ws.Range("S39").Value = 1293
$ws.Range("T39").Value = 15

# Row 40 (was row 35)
$ws.Range("D40").Value = 44676
$ws.Range("K40").Value = 'Fuyu'
$ws.Range("L40").Value = 'Primera'
$ws.Range("M40").Value = 115
$ws.Range("N40").Value = 15000
$ws.Range("O40").Value = 15000
$ws.Range("P40").Value = 15000
$ws.Range("Q40").Value = '$/bandeja 15 kilos'
$ws.Range("R40").Value = 'Región de O''Higgins'
$ws.Range("S40").Value = 1000
$ws.Range("T40").Value = 15

# Row 41 (was row 26)
$ws.Range("D41").Value = 44313
$ws.Range("K41").Value = 'Fuyu'
$ws.Range("L41").Value = 'Primera'
$ws.Range("M41").Value = 35
$ws.Range("N41").Value = 15000
$ws.Range("O41").Value = 15000
$ws.Range("P41").Value = 15000
$ws.Range("Q41").Value = '$/bandeja 15 kilos granel'
$ws.Range("R41").Value = 'Provincia de Limarí'
$ws.Range("S41").Value = 1000
$ws.Range("T41").Value = 15

# Row 43 (was row 41)
$ws.Range("D43").Value = 44327
$ws.Range("K43").Value = 'Mankaki'
$ws.Range("L43").Value = 'Primera'
$ws.Range("M43").Value = 35
$ws.Range("N43").Value = 20000
$ws.Range("O43").Value = 20000
$ws.Range("P43").Value = 20000
$ws.Range("Q43").Value = '$/bandeja 15 kilos granel'
$ws.Range("R43").Value = 'Región de O''Higgins'
$ws.Range("S43").Value = 1333
$ws.Range("T43").Value = 15

# Row 44 (was row 36)
$ws.Range("D44").Value = 44348
$ws.Range("K44").Value = 'Hachiya'
$ws.Range("L44").Value = 'Primera'
$ws.Range("M44").Value = 30
$ws.Range("N44").Value = 25000
$ws.Range("O44").Value = 25000
$ws.Range("P44").Value = 25000
$ws.Range("Q44").Value = '$/caja 18 kilos granel'
$ws.Range("R44").Value = 'Región de O''Higgins'
$ws.Range("S44").Value = 1389
$ws.Range("T44").Value = 18

# Row 45 (was row 37)
$ws.Range("D45").Value = 44348
$ws.Range("K45").Value = 'Mankaki'
$ws.Range("L45").Value = 'Primera'
$ws.Range("M45").Value = 40
$ws.Range("N45").Value = 25000
$ws.Range("O45").Value = 25000
$ws.Range("P45").Value = 25000
$ws.Range("Q45").Value = '$/caja 18 kilos granel'
$ws.Range("R45").Value = 'Región de O''Higgins'
$ws.Range("S45").Value = 1389
$ws.Range("T45").Value = 18

# Row 46 (was row 27)
$ws.Range("D46").Value = 44341
$ws.Range("K46").Value = 'Mankaki'
$ws.Range("L46").Value = 'Primera'
$ws.Range("M46").Value = 25
$ws.Range("N46").Value = 17000
$ws.Range("O46").Value = 17000
$ws.Range("P46").Value = 17000
$ws.Range("Q46").Value = '$/bandeja 15 kilos granel'
$ws.Range("R46").Value = 'Provincia de Limarí'
$ws.Range("S46").Value = 1133
$ws.Range("T46").Value = 15

# Row 47 (was row 19)
$ws.Range("D47").Value = 44708
$ws.Range("K47").Value = 'Mankaki'
$ws.Range("L47").Value = 'Primera'
$ws.Range("M47").Value = 45
$ws.Range("N47").Value = 17000
$ws.Range("O47").Value = 18000
$ws.Range("P47").Value = 17444
$ws.Range("Q47").Value = '$/caja 18 kilos granel'
$ws.Range("R47").Value = 'Región de O''Higgins'
$ws.Range("S47").Value = 969
$ws.Range("T47").Value = 18

# Row 48 (was row 34)
$ws.Range("D48").Value = 44301
$ws.Range("K48").Value = 'Fuyu'
$ws.Range("L48").Value = 'Primera'
$ws.Range("M48").Value = 75
$ws.Range("N48").Value = 10000
$ws.Range("O48").Value = 10000
$ws.Range("P48").Value = 10000
$ws.Range("Q48").Value = '$/caja 10 kilos granel'
$ws.Range("R48").Value = 'Región de O''Higgins'
$ws.Range("S48").Value = 10000
$ws.Range("T48").Value = 1

# Row 49 (was row 28)
$ws.Range("D49").Value = 44336
$ws.Range("K49").Value = 'Fuyu'
$ws.Range("L49").Value = 'Especial'
$ws.Range("M49").Value = 30
$ws.Range("N49").Value = 18000
$ws.Range("O49").Value = 18000
$ws.Range("P49").Value = 18000
$ws.Range("Q49").Value = '$/caja 15 kilos empedrada'
$ws.Range("R49").Value = 'Provincia de Limarí'
$ws.Range("S49").Value = 1200
$ws.Range("T49").Value = 15

# Row 50 (was row 29)
$ws.Range("D50").Value = 44336
$ws.Range("K50").Value = 'Mankaki'
$ws.Range("L50").Value = 'Especial'
$ws.Range("M50").Value = 50
$ws.Range("N50").Value = 18000
$ws.Range("O50").Value = 18000
$ws.Range("P50").Value = 18000
$ws.Range("Q50").Value = '$/caja 15 kilos empedrada'
$ws.Range("R50").Value = 'Provincia de Limarí'
$ws.Range("S50").Value = 1200
$ws.Range("T50").Value = 15

# Row 51 (was row 40)
$ws.Range("D51").Value = 44298
$ws.Range("K51").Value = 'Fuyu'
$ws.Range("L51").Value = 'Primera'
$ws.Range("M51").Value = 95
$ws.Range("N51").Value = 10000
$ws.Range("O51").Value = 10000
$ws.Range("P51").Value = 10000
$ws.Range("Q51").Value = '$/bandeja 15 kilos granel'
$ws.Range("R51").Value = 'Región de O''Higgins'
$ws.Range("S51").Value = 667
$ws.Range("T51").Value = 15

# Row 52 (was row 51)
$ws.Range("D52").Value = 44340
$ws.Range("K52").Value = 'Mankaki'
$ws.Range("L52").Value = 'Primera'
$ws.Range("M52").Value = 85
$ws.Range("N52").Value = 16000
$ws.Range("O52").Value = 17000
$ws.Range("P52").Value = 16471
$ws.Range("Q52").Value = '$/bandeja 15 kilos granel'
$ws.Range("R52").Value = 'Provincia de Limarí'
$ws.Range("S52").Value = 1098
$ws.Range("T52").Value = 15
